{"js": "// Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n// (percentages, dollar amounts, large numbers) in specific bullet points,\n// matching the target diff exactly.\n//\n// Strategy: for each target paragraph (identified by its exact original\n// text, to avoid any substring ambiguity between similar bullets), search\n// *within that paragraph* for each metric token in order and set its font\n// to bold + color 2C3E50. Word/Office.js automatically splits the run\n// containing the hit into separate runs, leaving the surrounding plain\n// text in its own run(s) \u2014 the same run-splitting behavior the target\n// diff shows.\n\nconst HILITE_COLOR = \"#2C3E50\";\n\nasync function highlightInParagraph(context, paragraph, tokens) {\n  // Re-search (and re-sync) per token so earlier splits in this paragraph\n  // don't invalidate later range objects.\n  for (const token of tokens) {\n    const results = paragraph.search(token, { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n    if (results.items.length === 0) {\n      continue;\n    }\n    // Each token is unique (single occurrence) within its paragraph in this\n    // document, so the first hit is the one we want.\n    const hit = results.items[0];\n    hit.font.set({ bold: true, color: HILITE_COLOR });\n    await context.sync();\n  }\n}\n\n// Each entry: the exact original text of the target paragraph, plus the\n// ordered list of metric substrings to bold+color within it.\nconst targets = [\n  {\n    text: \"\\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n    tokens: [\"23%\", \"64%\"],\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \\u00B14.2% to \\u00B12.1%\",\n    tokens: [\"87%\", \"71%\", \"\\u00B14.2%\", \"\\u00B12.1%\"],\n  },\n  {\n    text: \"\\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n    tokens: [\"1,200\"],\n  },\n  {\n    text: \"\\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n    tokens: [\"$400M\", \"$1B\"],\n  },\n  {\n    text: \"\\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n    tokens: [\"73.5%\", \"$4.7M\"],\n  },\n  {\n    text: \"\\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n    tokens: [\"87%\", \"71%\"],\n  },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const target of targets) {\n  const paragraph = paragraphs.items.find((p) => p.text === target.text);\n  if (!paragraph) {\n    continue;\n  }\n  await highlightInParagraph(context, paragraph, target.tokens);\n}\n", "ps1": "# Apply \"hybrid bold + color\" highlighting to quantitative impact metrics\n# (percentages, dollar amounts, large numbers) in specific bullet points,\n# matching the target diff exactly.\n#\n# Strategy: for each target paragraph (identified by its exact original\n# text, to avoid any substring ambiguity between similar bullets), use\n# Find.Execute *scoped to that paragraph's Range* to locate each metric\n# token in turn and set its Font.Bold/Font.Color. Word splits the run\n# containing the hit into separate runs, leaving the surrounding plain\n# text in its own run(s) - the same run-splitting behavior the target\n# diff shows.\n\nfunction Get-WdColor($r, $g, $b) {\n    # Word's Font.Color (wdColor) packs RGB as 0x00BBGGRR.\n    return $b * 65536 + $g * 256 + $r\n}\n\n$HiliteColor = Get-WdColor 0x2C 0x3E 0x50\n\nfunction Set-MetricHighlights($paragraphRange, $tokens) {\n    foreach ($token in $tokens) {\n        # Re-duplicate the paragraph range fresh for every search so the\n        # find is scoped correctly and isn't thrown off by earlier splits.\n        $findRange = $paragraphRange.Duplicate\n        $found = $findRange.Find.Execute($token)\n        if ($found) {\n            $findRange.Font.Bold = $true\n            $findRange.Font.Color = $HiliteColor\n        }\n    }\n}\n\n# Each entry: the exact original text of the target paragraph (paragraph\n# mark excluded), plus the ordered list of metric substrings to\n# bold+color within it.\n$targets = @(\n    @{\n        Text   = \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\"\n        Tokens = @(\"23%\", \"64%\")\n    },\n    @{\n        Text   = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\"\n        Tokens = @(\"87%\", \"71%\", \"\u00b14.2%\", \"\u00b12.1%\")\n    },\n    @{\n        Text   = \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\"\n        Tokens = @(\"1,200\")\n    },\n    @{\n        Text   = \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\"\n        Tokens = @(\"`$400M\", \"`$1B\")\n    },\n    @{\n        Text   = \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\"\n        Tokens = @(\"73.5%\", \"`$4.7M\")\n    },\n    @{\n        Text   = \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n        Tokens = @(\"87%\", \"71%\")\n    }\n)\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\nforeach ($target in $targets) {\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $t = $p.Range.Text.TrimEnd([char]13)\n        if ($t -eq $target.Text) {\n            Set-MetricHighlights $p.Range $target.Tokens\n            break\n        }\n    }\n}\n"}
